$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 820.5217
$ws.Range("I33").Value = 287.3889
$ws.Range("K33").Value = 287.3889
$ws.Range("M33").Value = -58.38889999999998
$ws.Range("H64").Value = 2627.5
$ws.Range("I64").Value = 2514
$ws.Range("K64").Value = 2514
$ws.Range("M64").Value = -2266
$ws.Range("H67").Value = 2627.5
$ws.Range("I67").Value = 2514
$ws.Range("K67").Value = 2514
$ws.Range("M67").Value = -1656
$ws.Range("H116").Value = 7793.4585
$ws.Range("I116").Value = 6773.077
$ws.Range("K116").Value = 6773.077
$ws.Range("M116").Value = -3331.077
$ws.Range("H137").Value = 1184.25
$ws.Range("I137").Value = 1168
$ws.Range("K137").Value = 3504
$ws.Range("M137").Value = -954

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 31999.445
$ws.Range("I33").Value = 19655.666
$ws.Range("K33").Value = 19655.666
$ws.Range("M33").Value = -19326.666
$ws.Range("H61").Value = 4586.143
$ws.Range("I61").Value = 4985
$ws.Range("K61").Value = 4985
$ws.Range("M61").Value = -4773
$ws.Range("H74").Value = 1602
$ws.Range("I74").Value = 1557.7778
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 1557.7778
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -683.7778000000001
$ws.Range("N74").Value = -3748
$ws.Range("H77").Value = 1602
$ws.Range("I77").Value = 1557.7778
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 7788.889
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -3420.889
$ws.Range("N77").Value = -18736
$ws.Range("H122").Value = 2147.1738
$ws.Range("I122").Value = 2252.0527
$ws.Range("K122").Value = 6756.158100000001
$ws.Range("M122").Value = -4306.158100000001
$ws.Range("H132").Value = 2221.1233
$ws.Range("I132").Value = 2190
$ws.Range("K132").Value = 6570
$ws.Range("M132").Value = -4040
$ws.Range("H136").Value = 4586.143
$ws.Range("I136").Value = 4985
$ws.Range("K136").Value = 14955
$ws.Range("M136").Value = -12405

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H99").Value = 1929.4445
$ws.Range("I99").Value = 1974.2354
$ws.Range("K99").Value = 1974.2354
$ws.Range("M99").Value = -476.2354
$ws.Range("H107").Value = 12114.366
$ws.Range("I107").Value = 13775.68
$ws.Range("J107").Value = 3807.8
$ws.Range("K107").Value = 13775.68
$ws.Range("L107").Value = 3807.8
$ws.Range("M107").Value = -11855.68
$ws.Range("N107").Value = -7647.8
$ws.Range("H134").Value = 4149.8
$ws.Range("I134").Value = 3824.5625
$ws.Range("K134").Value = 11473.6875
$ws.Range("M134").Value = -8938.6875
$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value = 99899.5
$ws.Range("J137").Value = 99899.5
$ws.Range("L137").Value = 99899.5
$ws.Range("N137").Value = -110099.5
$ws.Range("H140").Value = 128572.7
$ws.Range("J140").Value = 128572.7
$ws.Range("L140").Value = 128572.7
$ws.Range("N140").Value = -138932.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1754.1818
$ws.Range("I31").Value = 868.375
$ws.Range("J31").Value = 3304.3438
$ws.Range("K31").Value = 868.375
$ws.Range("L31").Value = 3304.3438
$ws.Range("M31").Value = -573.375
$ws.Range("N31").Value = -3894.3438
$ws.Range("H34").Value = 1754.1818
$ws.Range("I34").Value = 868.375
$ws.Range("J34").Value = 3304.3438
$ws.Range("K34").Value = 868.375
$ws.Range("L34").Value = 3304.3438
$ws.Range("M34").Value = -666.375
$ws.Range("N34").Value = -3708.3438
$ws.Range("H62").Value = 2960.7144
$ws.Range("I62").Value = 2833.75
$ws.Range("J62").Value = 3130
$ws.Range("K62").Value = 2833.75
$ws.Range("L62").Value = 3130
$ws.Range("M62").Value = -2209.75
$ws.Range("N62").Value = -4378
$ws.Range("H65").Value = 2960.7144
$ws.Range("I65").Value = 2833.75
$ws.Range("J65").Value = 3130
$ws.Range("K65").Value = 14168.75
$ws.Range("L65").Value = 15650
$ws.Range("M65").Value = -11048.75
$ws.Range("N65").Value = -21890
$ws.Range("H86").Value = 9741.691999999999
$ws.Range("I86").Value = 5908
$ws.Range("K86").Value = 5908
$ws.Range("M86").Value = -4785
$ws.Range("H89").Value = 9741.691999999999
$ws.Range("I89").Value = 5908
$ws.Range("K89").Value = 29540
$ws.Range("M89").Value = -23924
$ws.Range("H99").Value = 4095.9167
$ws.Range("I99").Value = 3577.7144
$ws.Range("J99").Value = 4821.4
$ws.Range("K99").Value = 3577.7144
$ws.Range("L99").Value = 4821.4
$ws.Range("M99").Value = -2079.7144
$ws.Range("N99").Value = -7817.4
$ws.Range("H122").Value = 3327.9375
$ws.Range("I122").Value = 3327.9375
$ws.Range("K122").Value = 9983.8125
$ws.Range("M122").Value = -7533.8125
$ws.Range("H126").Value = 4095.9167
$ws.Range("I126").Value = 3577.7144
$ws.Range("J126").Value = 4821.4
$ws.Range("K126").Value = 10733.1432
$ws.Range("L126").Value = 14464.2
$ws.Range("M126").Value = -8263.143199999999
$ws.Range("N126").Value = -19404.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 4992.2085
$ws.Range("I134").Value = 4753.316
$ws.Range("K134").Value = 14259.948
$ws.Range("M134").Value = -9189.948
$ws.Range("H140").Value = 1032.4286
$ws.Range("I140").Value = 1032.4286
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 3097.2858
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = 2082.7142
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 59211.69
$ws.Range("J64").Value = 59211.69
$ws.Range("L64").Value = 59211.69
$ws.Range("N64").Value = -59707.69
$ws.Range("H67").Value = 59211.69
$ws.Range("J67").Value = 59211.69
$ws.Range("L67").Value = 59211.69
$ws.Range("N67").Value = -60927.69
$ws.Range("H102").Value = 5620
$ws.Range("I102").Value = 5764.7
$ws.Range("K102").Value = 5764.7
$ws.Range("M102").Value = -4142.7
$ws.Range("H126").Value = 4935.4443
$ws.Range("I126").Value = 4935.4443
$ws.Range("K126").Value = 14806.3329
$ws.Range("M126").Value = -12336.3329
$ws.Range("H132").Value = 940184.75
$ws.Range("I132").Value = 1156334.6
$ws.Range("K132").Value = 3469003.8
$ws.Range("M132").Value = -3466473.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 16787.666
$ws.Range("I7").Value = 17764.371
$ws.Range("K7").Value = 17764.371
$ws.Range("M7").Value = -17652.371
$ws.Range("H122").Value = 3600.8125
$ws.Range("I122").Value = 3367.5833
$ws.Range("K122").Value = 10102.7499
$ws.Range("M122").Value = -7652.749899999999
$ws.Range("H126").Value = 16787.666
$ws.Range("I126").Value = 17764.371
$ws.Range("K126").Value = 53293.113
$ws.Range("M126").Value = -50823.113
$ws.Range("H132").Value = 2462391.8
$ws.Range("I132").Value = 2845021.8
$ws.Range("K132").Value = 8535065.399999999
$ws.Range("M132").Value = -8532535.399999999
$ws.Range("H136").Value = 17549626
$ws.Range("I136").Value = 55561604
$ws.Range("J136").Value = 5634.3076
$ws.Range("K136").Value = 166684812
$ws.Range("L136").Value = 16902.9228
$ws.Range("M136").Value = -166682262
$ws.Range("N136").Value = -22002.9228

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4057.4417
$ws.Range("I122").Value = 3509.8225
$ws.Range("K122").Value = 10529.4675
$ws.Range("M122").Value = -8079.467500000001
$ws.Range("H136").Value = 15388152
$ws.Range("I136").Value = 19233940
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 57701820
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -57699270
$ws.Range("N136").Value = -20100
